# Slide 1, shape "Прямоугольник 5" (the subtitle rectangle), 4th paragraph:
#   "РОБОТИЗИРОВАННЫЙ " + " " + "МНОГОЦЕЛЕВОЙ КОМПЛЕКС" + " " + "ДЛЯ ПРОМЫШЛЕННЫХ ..."
#
# The diff removes the lone-space run (lang="en-US") that sits between
# "РОБОТИЗИРОВАННЫЙ " and "МНОГОЦЕЛЕВОЙ КОМПЛЕКС", and splits the
# "МНОГОЦЕЛЕВОЙ КОМПЛЕКС" run (lang="ru-RU") into two runs with identical
# rPr: "МНОГОЦЕЛЕВОЙ " and "КОМПЛЕКС".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Isolate the 4th paragraph as its own TextRange (1-based index, count=1).
$para4 = $tr.Paragraphs(4, 1)

# Within this paragraph (1-based, local offsets):
#   chars 1..17  -> "РОБОТИЗИРОВАННЫЙ "
#   char  18     -> " "                       (lone space run, to be removed)
#   chars 19..39 -> "МНОГОЦЕЛЕВОЙ КОМПЛЕКС"    (to be split 13/8)
#   char  40     -> " "
#   ...

# Step 1: split "МНОГОЦЕЛЕВОЙ КОМПЛЕКС" into "МНОГОЦЕЛЕВОЙ " + "КОМПЛЕКС".
# Touching identical formatting on the two halves forces PowerPoint to
# materialize them as two distinct runs sharing the same rPr, mirroring
# what the target OOXML shows.
$part1 = $para4.Characters(19, 13)   # "МНОГОЦЕЛЕВОЙ "
$part2 = $para4.Characters(32, 8)    # "КОМПЛЕКС"
$part1.Font.Name = "+mj-lt"
$part2.Font.Name = "+mj-lt"

# Step 2: remove the lone-space run between "РОБОТИЗИРОВАННЫЙ " and
# "МНОГОЦЕЛЕВОЙ" (still at local offset 18 - the split above only touched
# characters at/after offset 19, so this offset is unaffected).
$spaceRun = $para4.Characters(18, 1)
$spaceRun.Text = ""
